# Applies the weekly-update edit to the "Hortaliza, Vega Modelo de Temuco -
# Ciboulette" sheet:
#   - a brand-new price observation is inserted as row 137;
#   - every existing observation in rows 137..269 shifts down by one row
#     (row R takes what used to live in row R-1), for the columns that vary
#     per-observation: D (Fecha), I (Calidad), J (Volumen), K (Precio
#     minimo), L (Precio maximo), M (Precio promedio ponderado), O (Origen),
#     P (Precio $/Kg);
#   - the observation that used to be the last one (old row 269) overflows
#     into a brand-new row 270, keeping its constant/descriptive columns
#     (A,B,C,E,F,G,H,N,Q,R) the same as the rest of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 137
$lastRow  = 269
$newRow   = 270

# Columns (1-based) whose value is tied to a specific observation and thus
# shifts down by one row each week.
$shiftCols = @(4, 9, 10, 11, 12, 13, 15, 16)

# Columns that stay constant for every row in this block (copied verbatim
# into the newly appended row).
$constCols = @(1, 2, 3, 5, 6, 7, 8, 14, 17, 18)

# Remember the last row's per-observation values and number format before
# they get overwritten -- they become the new row 270.
$overflow = @{}
foreach ($c in $shiftCols) {
    $overflow[$c] = $ws.Cells.Item($lastRow, $c).Value2
}
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# Shift rows 138..269 down from 137..268, processing bottom-up so each
# source cell is read before it gets overwritten.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    foreach ($c in $shiftCols) {
        $prevVal = $ws.Cells.Item($r - 1, $c).Value2
        $ws.Cells.Item($r, $c).Value = $prevVal
    }
}

# Row 137 becomes the brand-new observation.
$ws.Cells.Item($firstRow, 4).Value  = 44778
$ws.Cells.Item($firstRow, 9).Value  = "Primera"
$ws.Cells.Item($firstRow, 10).Value = 45
$ws.Cells.Item($firstRow, 11).Value = 7000
$ws.Cells.Item($firstRow, 12).Value = 7000
$ws.Cells.Item($firstRow, 13).Value = 7000
$ws.Cells.Item($firstRow, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($firstRow, 16).Value = 2333

# Append row 270: constant columns copied from row 269, per-observation
# columns filled with the captured overflow values.
foreach ($c in $constCols) {
    $ws.Cells.Item($newRow, $c).Value = $ws.Cells.Item($lastRow, $c).Value2
}
foreach ($c in $shiftCols) {
    $ws.Cells.Item($newRow, $c).Value = $overflow[$c]
}
$ws.Cells.Item($newRow, 4).NumberFormat = $dateFormat
